$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells I8:I18 already use the numeric style that the new J column cells
# need (numFmt "0.0", Times New Roman 9pt, vertically centered, no border).
# Reuse that formatting for the newly populated J cells (rows 19-26).
$target = $ws.Range("J19:J26")
$target.NumberFormat = "0.0"
$target.Font.Name = "Times New Roman"
$target.Font.Size = 9
$target.VerticalAlignment = -4108

# Fill in new column J values for rows 19-27
$ws.Range("J19").Value = 12.434613462352335
$ws.Range("J20").Value = 16.80050595536094
$ws.Range("J21").Value = 11.282963378125267
$ws.Range("J22").Value = 25.042808754677555
$ws.Range("J23").Value = 3.2011163356916352
$ws.Range("J24").Value = 13.523574517571838
$ws.Range("J25").Value = 6.1196997869329204
$ws.Range("J26").Value = 5.9488136666578013

# J27 already has the correct style (s="30"); just populate its value.
$ws.Range("J27").Value = 5.2451982064110645

# Update the active selection to N8
$ws.Range("N8").Select()
